# Toggle the "TO BE EXECUTED" (column A) YES/NO flags for the affected
# test-step rows, then update the sheet's zoom level and active selection
# to match the reviewer's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose column-A flag flips from YES -> NO
$rowsToNo = @(2, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25)

# Rows whose column-A flag flips from NO -> YES
$rowsToYes = @(38, 39, 40, 41, 42, 43, 44)

foreach ($r in $rowsToNo) {
    $ws.Range("A$r").Value = "NO"
}

foreach ($r in $rowsToYes) {
    $ws.Range("A$r").Value = "YES"
}

# Update view state: zoom level and active cell/selection
$ws.Application.ActiveWindow.Zoom = 79
$ws.Range("A2").Select() | Out-Null
